$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Add a new worksheet "Sheet, with comma" right after "Global".
# Its content duplicates the "Global" sheet's header row + first data row
# (same Container Path / Parameter Name / Value / Units layout).
# ---------------------------------------------------------------------------
$global = $wb.Worksheets.Item("Global")
$newSheet = $wb.Worksheets.Add($null, $global)
$newSheet.Name = "Sheet, with comma"

$newSheet.Range("A1").Value = "Container Path"
$newSheet.Range("B1").Value = "Parameter Name"
$newSheet.Range("C1").Value = "Value"
$newSheet.Range("D1").Value = "Units"
$newSheet.Range("A2").Value = "Organism|Liver"
$newSheet.Range("B2").Value = "EHC continuous fraction"
$newSheet.Range("C2").Value = 1
$newSheet.Range("C1:C2").NumberFormat = "0.0000"

# Match the column layout used on the "Global" sheet it was copied from.
$newSheet.Columns.Item(1).ColumnWidth = $global.Columns.Item(1).ColumnWidth
$newSheet.Columns.Item(2).ColumnWidth = $global.Columns.Item(2).ColumnWidth
$newSheet.Columns.Item(3).ColumnWidth = $global.Columns.Item(3).ColumnWidth
$newSheet.Columns.Item(4).ColumnWidth = $global.Columns.Item(4).ColumnWidth

# ---------------------------------------------------------------------------
# The new sheet becomes the active tab; the view is scrolled to A30.
# ---------------------------------------------------------------------------
[void]$newSheet.Activate()
[void]$newSheet.Range("A30").Select()

# ---------------------------------------------------------------------------
# "Global" is no longer the selected tab; its selection now spans every row.
# ---------------------------------------------------------------------------
[void]$global.Activate()
[void]$global.Cells.Select()

# Leave the new sheet active/selected, matching the saved workbook view.
[void]$newSheet.Activate()
